# Q3 Update - 2025
# Refresh the UN-SOM (Somalia) refugee-statistics export with the latest
# Q3 2025 figures: a new dataset "short-url" token for every data row plus
# a handful of updated statistic cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- short-url column (B) refresh for every data row (2..305) ---
for ($r = 2; $r -le 305; $r++) {
    $ws.Cells.Item($r, 2).Value = "5lQrA3"
}

# --- statistic cell updates -------------------------------------------------
# These columns hold numeric-looking values that are stored as text in the
# source data; force text format first so Excel keeps writing them as text
# (matching the existing column formatting) instead of auto-converting them
# to numbers.
$statCells = @("O291","O294","O295","N297","R301","T301","N302","O302","N303","O303","N305")
foreach ($addr in $statCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 291 - Egypt -> Somalia : asylum_seekers 5 -> 6
$ws.Range("O291").Value = "6"

# Row 294 - Djibouti -> Somalia : asylum_seekers 5 -> 6
$ws.Range("O294").Value = "6"

# Row 295 - Eritrea -> Somalia : asylum_seekers 53 -> 54
$ws.Range("O295").Value = "54"

# Row 297 - Palestinian -> Somalia : refugees 9 -> 12
$ws.Range("N297").Value = "12"

# Row 301 - Somalia -> Somalia : returned_idps 0 -> 290750, ooc 222 -> 213
$ws.Range("R301").Value = "290750"
$ws.Range("T301").Value = "213"

# Row 302 - Sudan -> Somalia : refugees 214 -> 273, asylum_seekers 25 -> 42
$ws.Range("N302").Value = "273"
$ws.Range("O302").Value = "42"

# Row 303 - Syrian Arab Rep. -> Somalia : refugees 1545 -> 1591, asylum_seekers 164 -> 178
$ws.Range("N303").Value = "1591"
$ws.Range("O303").Value = "178"

# Row 305 - Yemen -> Somalia : refugees 11647 -> 12332
$ws.Range("N305").Value = "12332"
